$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "6.915 - x - 0.5y"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "-4.915"
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "0.8999999999999999"
$ws.Range("F2").Style = "Normal"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "-1.0500000000000003 - 0.25x + y"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "-0.9499999999999997"
$ws.Range("B3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.4000000000000004"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "0.8999999999999999"
$ws.Range("F3").Style = "Normal"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "-6.915 + x + 0.5y"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "-1.085"
$ws.Range("B4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.97"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "3.9000000000000004"
$ws.Range("F4").Style = "Normal"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "-3.2600000000000007 + x - 2y"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "-1.2600000000000007"
$ws.Range("B5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.2"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "9.399999999999999"
$ws.Range("F5").Style = "Normal"
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "-2.47 - y"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "-2.47"
$ws.Range("B6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "8.5"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2.1"
$ws.Range("F6").Style = "Normal"

$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "5.68"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2.47"
$ws.Range("B2").Style = "Normal"

$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1.76"
$ws.Range("A2").Style = "Normal"

$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "-1.3499999999999999"
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "6.5"
$ws.Range("A3").Style = "Normal"
